# Update the cryptos list (Sat Oct 14 19:16:16 UTC 2023 GitHub Actions refresh).
# Updates Price (D) and Volume(1h) (E) values for each coin row, and swaps the
# Maker / InternetComputer(DFINITY) rows (33 <-> 34).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.974.62'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '1.560.18'
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.38'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.489'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.11'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.80%  '
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("E10").Value = '  +1.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0859'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '1.783.01'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '1.545.63'
$ws.Range("E13").Value = '  -0.56%  '
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("D17").Value = '26.977.53'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.20'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").Value = '0.0₃0703'
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.37'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.21%  '
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.09'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.20'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.63'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.61'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("E27").Value = '  +1.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.104'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("E31").Value = '  +2.01%  '
$ws.Range("E32").Value = '  +0.65%  '
$ws.Range("B33").Value = 'Maker'
$ws.Range("C33").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D33").Value = '1.422.44'
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.11'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.60%  '
$ws.Range("E35").Value = '  +2.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +9.26%  '
$ws.Range("E37").Value = '  +1.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0165'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.534'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.808'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.70'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("E43").Value = '  +2.97%  '
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.85'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.86%  '
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").Value = '1.696.09'
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.39'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0522'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0953'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.40%  '
